# Align WHO Rotavirus workbook with current WHO guidance:
# Remove Dose 1 maxAge "15 weeks" restriction (WHO 2013 relaxation of age
# restrictions) and qualify the series-name sheet tabs with their vaccine
# abbreviations (RV1 / RV5). The 8-month hard ceiling on later doses is
# left untouched.

$wb = $excel.ActiveWorkbook

# --- 1) Rename the two series tabs ---------------------------------------
$ws2 = $wb.Worksheets.Item("2-dose series")
$ws2.Name = "2-dose series (RV1)"

$ws3 = $wb.Worksheets.Item("3-dose series")
$ws3.Name = "3-dose series (RV5)"

# --- 2) "2-dose series (RV1)" sheet: Dose 1 Age row (row 8) --------------
# Column F held the Dose-1 max age ("15 weeks"); WHO no longer restricts
# this, so it becomes "n/a" like the other unused bound columns, and the
# now-redundant trailing n/a cells (G8:H8) are dropped entirely.
$ws2.Range("F8").Value = "n/a"
$ws2.Range("G8:H8").Clear()

# Trim the redundant trailing "n/a" cells that trailed the Dose 2 Age /
# Preferable Interval rows (content unchanged, just fewer duplicate cells).
$ws2.Range("G14:H14").Clear()
$ws2.Range("J15:L15").Clear()

# --- 3) "3-dose series (RV5)" sheet: same Dose 1 Age fix (row 8) ---------
$ws3.Range("F8").Value = "n/a"
$ws3.Range("G8:H8").Clear()

# Trim the same redundant trailing "n/a" cells for Dose 2 ...
$ws3.Range("G14:H14").Clear()
$ws3.Range("J15:L15").Clear()

# ... and for Dose 3 (rows 21/22 on this sheet only).
$ws3.Range("G21:H21").Clear()
$ws3.Range("J22:L22").Clear()
